$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New layout: A Month | B Date | C Employee ID | D Employee Name |
# E Check In | F Check In Status | G Check In Time |
# H Check Out | I Check Out Status | J Check Out Time | K Total Salary
#
# NOTE: new shared-string values below are introduced in a specific
# order (Check In Time, Check Out Time, 11/7/2023, 8:54, 19:55, late,
# 10:55, No, missing, 21:56) so the rebuilt shared string table lines
# up with the target workbook.

# --- Header row (existing headers reuse old shared strings; only
#     "Check In Time" and "Check Out Time" are new) ---
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Employee ID"
$ws.Range("D1").Value = "Employee Name"
$ws.Range("E1").Value = "Check In"
$ws.Range("F1").Value = "Check In Status"
$ws.Range("G1").Value = "Check In Time"
$ws.Range("H1").Value = "Check Out"
$ws.Range("I1").Value = "Check Out Status"
$ws.Range("J1").Value = "Check Out Time"
$ws.Range("K1").Value = "Total Salary"

# --- Row 2 (ThaiNX) ---
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "'11/7/2023"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "ThaiNX"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "on time"
$ws.Range("G2").Value = "8:54"
$ws.Range("H2").Value = "Yes"
$ws.Range("I2").Value = "on time"
$ws.Range("J2").Value = "19:55"
$ws.Range("K2").Value = 1000000

# --- Row 3 (Thanh Ha) ---
$ws.Range("D3").Value = "Thanh Ha"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "late"
$ws.Range("G3").Value = "10:55"
$ws.Range("H3").Value = "No"
$ws.Range("I3").Value = "missing"
$ws.Range("J3").Value = "21:56"
$ws.Range("A3").Clear()
$ws.Range("B3").Clear()
$ws.Range("C3").Value = 2
$ws.Range("K3").Value = 500000

# --- Row 4 (Duy Long) removed entirely ---
$ws.Range("A4:K4").Clear()

# Extend widened columns E:I (width 15) to cover new columns through K.
# (14.1666... round-trips to a stored column width of exactly 15, matching
# the width already used for columns E:I.)
$ws.Range("E1:K1").ColumnWidth = 14.1666666666667
